$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 8 ("franzosa_ControlvsCD_ConvCD"), pushing
# "franzosa_ControlvsCD_Fp" and everything below it down by one row.
$ws.Rows.Item(8).Insert()

# Insert a new row at row 13 ("franzosa_ControlvsUC_ConvUC"). At this
# point "franzosa_ControlvsUC_Age" has already shifted down to row 12
# because of the first insert above, so this pushes "yachida_age" and
# everything below it down by one more row.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 8.
$ws.Range("A8").Value = "franzosa_ControlvsCD_ConvCD"
$ws.Range("B8").Value = 0.05
$ws.Range("C8").Value = 0.05
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0.55
$ws.Range("F8").Value = 0.9
$ws.Range("G8").Value = 0.35
$ws.Range("H8").Value = 0.4

# Populate the newly inserted row 13.
$ws.Range("A13").Value = "franzosa_ControlvsUC_ConvUC"
$ws.Range("B13").Value = 0.05
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0.75
$ws.Range("F13").Value = 0.95
$ws.Range("G13").Value = 0.2
$ws.Range("H13").Value = 0.2
